# Update the AUGUST_2024 attendance sheet: several students' attendance
# marks are corrected from 0 to 3, which ripples through the SUM()
# (column G) and percentage (column H) formulas automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AUGUST_2024")
$ws.Activate()

# Roll 23ME01 (row 14): E14 0 -> 3
$ws.Range("E14").Value = 3

# Roll 23ME04 (row 17): D17 0 -> 3, E17 0 -> 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 3

# Roll 23ME06 (row 19): E19 0 -> 3
$ws.Range("E19").Value = 3

# Roll 23ME07 (row 20): E20 0 -> 3
$ws.Range("E20").Value = 3

# Roll 23ME09 (row 21): D21 0 -> 3
$ws.Range("D21").Value = 3

# Roll 22ME03 (row 22): E22 0 -> 3
$ws.Range("E22").Value = 3

# Roll 21ME20 (row 24): E24 0 -> 3
$ws.Range("E24").Value = 3

# Leave the active selection where the editor finished up.
$ws.Range("D18").Select()
